$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellAddr, $val) {
    $rng = $ws.Range($cellAddr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue "D2" '309.88'
Set-TextValue "E2" '-3.02%'
Set-TextValue "D3" '54.22'
Set-TextValue "E3" '10.31%'
Set-TextValue "D4" '5.127'
Set-TextValue "E4" '-2.49%'
Set-TextValue "D5" '0.07831'
Set-TextValue "E5" '-1.55%'
Set-TextValue "D6" '4.523'
Set-TextValue "E6" '-1.09%'
Set-TextValue "D7" '1.364'
Set-TextValue "E7" '-3.02%'
Set-TextValue "D9" '0.1221'
Set-TextValue "E9" '-5.96%'
Set-TextValue "D10" '0.2011'
Set-TextValue "E10" '2.08%'
Set-TextValue "D11" '0.04717'
Set-TextValue "E11" '2.49%'
Set-TextValue "D12" '0.09458'
Set-TextValue "E12" '0.09%'
Set-TextValue "D13" '0.1045'
Set-TextValue "E13" '-0.20%'
Set-TextValue "D14" '0.001262'
Set-TextValue "E14" '-5.58%'
Set-TextValue "D15" '0.005770'
Set-TextValue "E15" '-2.86%'
Set-TextValue "E16" '2,018.14%'
Set-TextValue "D17" '3.335'
Set-TextValue "E17" '-0.25%'
Set-TextValue "D18" '2.414'
Set-TextValue "E18" '-0.87%'
Set-TextValue "D19" '0.3418'
Set-TextValue "E19" '-1.20%'
Set-TextValue "D20" '8.039'
Set-TextValue "E20" '-2.34%'
Set-TextValue "D21" '0.1369'
Set-TextValue "E21" '-1.47%'
Set-TextValue "D22" '0.3084'
Set-TextValue "E22" '-0.19%'
Set-TextValue "D23" '0.04166'
Set-TextValue "E23" '0.18%'
Set-TextValue "D24" '0.001259'
Set-TextValue "E24" '-4.24%'
Set-TextValue "D25" '0.003925'
Set-TextValue "E25" '-7.85%'
Set-TextValue "D26" '0.0001348'
Set-TextValue "E26" '-0.05%'
Set-TextValue "D38" '0.02607'
Set-TextValue "E38" '-3.02%'
Set-TextValue "D39" '0.05952'
Set-TextValue "E39" '2.79%'
Set-TextValue "E40" '-3.55%'
Set-TextValue "D41" '0.007916'
Set-TextValue "E41" '-1.09%'
Set-TextValue "D42" '0.1425'
Set-TextValue "E42" '-0.96%'
Set-TextValue "D43" '0.008224'
Set-TextValue "E43" '6.61%'
Set-TextValue "D44" '0.008478'
Set-TextValue "E44" '-0.13%'
Set-TextValue "D45" '0.3120'
Set-TextValue "E45" '-2.24%'
Set-TextValue "D46" '0.00007288'
Set-TextValue "E46" '10.28%'
Set-TextValue "D47" '0.00000000748'
Set-TextValue "E47" '-0.15%'
Set-TextValue "D48" '0.05632'
Set-TextValue "E48" '2.53%'
Set-TextValue "D49" '0.002613'
Set-TextValue "E49" '-34.59%'
Set-TextValue "D50" '0.00002095'
Set-TextValue "E50" '-0.15%'
Set-TextValue "D51" '0.0001995'
Set-TextValue "E51" '-0.15%'
